$wb = $excel.ActiveWorkbook

# --- Sheet "Erlaeuterung": update "Datenstand" and report date references ---
$wsInfo = $wb.Worksheets.Item("Erläuterung")
$wsInfo.Range("A6").Value = "Datenstand: 2021-04-26, 00:00 Uhr"
$wsInfo.Range("A11").Value = "Die Schätzwerte entsprechen der im täglichen Lagebericht des Robert Koch-Instituts vom 2021-04-26 dargestellten Nowcasting-Kurve."

# --- Sheet "Nowcast_R": revise nowcasting values for existing rows (RKI daily re-estimation) ---
$ws = $wb.Worksheets.Item("Nowcast_R")

$ws.Range("C356").Value = 7540
$ws.Range("C357").Value = 7281; $ws.Range("D357").Value = 7492
$ws.Range("B358").Value = 7125; $ws.Range("C358").Value = 7024; $ws.Range("D358").Value = 7233; $ws.Range("F358").Value = 7361
$ws.Range("B359").Value = 8956; $ws.Range("C359").Value = 8825; $ws.Range("D359").Value = 9068
$ws.Range("C360").Value = 8389; $ws.Range("D360").Value = 8602; $ws.Range("F360").Value = 7880; $ws.Range("G360").Value = 8098
$ws.Range("C361").Value = 7989; $ws.Range("D361").Value = 8205; $ws.Range("F361").Value = 8057; $ws.Range("G361").Value = 8277
$ws.Range("B362").Value = 7659; $ws.Range("D362").Value = 7771; $ws.Range("F362").Value = 8187; $ws.Range("G362").Value = 8411
$ws.Range("C363").Value = 7832; $ws.Range("D363").Value = 8061; $ws.Range("G363").Value = 8160
$ws.Range("B364").Value = 7181; $ws.Range("C364").Value = 7083; $ws.Range("D364").Value = 7282; $ws.Range("F364").Value = 7613
$ws.Range("B365").Value = 7344; $ws.Range("D365").Value = 7431; $ws.Range("F365").Value = 7422
$ws.Range("B366").Value = 9527; $ws.Range("C366").Value = 9408; $ws.Range("D366").Value = 9651
$ws.Range("B367").Value = 8912; $ws.Range("C367").Value = 8810; $ws.Range("D367").Value = 9038; $ws.Range("F367").Value = 8132; $ws.Range("G367").Value = 8350
$ws.Range("B368").Value = 8995; $ws.Range("C368").Value = 8863; $ws.Range("D368").Value = 9108; $ws.Range("G368").Value = 8807
$ws.Range("B369").Value = 8647; $ws.Range("C369").Value = 8531; $ws.Range("M369").Value = 1.07
$ws.Range("B370").Value = 9019; $ws.Range("D370").Value = 9152; $ws.Range("E370").Value = 8893; $ws.Range("F370").Value = 8779; $ws.Range("G370").Value = 9013
$ws.Range("B371").Value = 8901; $ws.Range("C371").Value = 8775; $ws.Range("D371").Value = 9025; $ws.Range("E371").Value = 8891; $ws.Range("F371").Value = 8771; $ws.Range("G371").Value = 9010
$ws.Range("B372").Value = 9032; $ws.Range("C372").Value = 8899; $ws.Range("D372").Value = 9157; $ws.Range("E372").Value = 8900; $ws.Range("F372").Value = 8780; $ws.Range("G372").Value = 9022
$ws.Range("B373").Value = 11793; $ws.Range("C373").Value = 11665; $ws.Range("D373").Value = 11934; $ws.Range("F373").Value = 9563; $ws.Range("G373").Value = 9817
$ws.Range("C374").Value = 11354; $ws.Range("D374").Value = 11633; $ws.Range("F374").Value = 10173
$ws.Range("B375").Value = 11730; $ws.Range("C375").Value = 11590; $ws.Range("D375").Value = 11863; $ws.Range("E375").Value = 11014; $ws.Range("F375").Value = 10877; $ws.Range("G375").Value = 11146
$ws.Range("B376").Value = 11497; $ws.Range("C376").Value = 11349; $ws.Range("D376").Value = 11647; $ws.Range("E376").Value = 11630; $ws.Range("F376").Value = 11489; $ws.Range("G376").Value = 11769; $ws.Range("I376").Value = 1.29
$ws.Range("B377").Value = 12059; $ws.Range("C377").Value = 11927; $ws.Range("D377").Value = 12203; $ws.Range("E377").Value = 11697; $ws.Range("F377").Value = 11555; $ws.Range("G377").Value = 11836
$ws.Range("C378").Value = 11694; $ws.Range("D378").Value = 11961; $ws.Range("E378").Value = 11779; $ws.Range("F378").Value = 11640; $ws.Range("G378").Value = 11918
$ws.Range("B379").Value = 12082; $ws.Range("C379").Value = 11962; $ws.Range("D379").Value = 12212; $ws.Range("E379").Value = 11867; $ws.Range("F379").Value = 11733; $ws.Range("G379").Value = 12006; $ws.Range("J379").Value = 1.08
$ws.Range("B380").Value = 15434; $ws.Range("C380").Value = 15285; $ws.Range("D380").Value = 15596; $ws.Range("E380").Value = 12851; $ws.Range("F380").Value = 12717; $ws.Range("G380").Value = 12993; $ws.Range("H380").Value = 1.1
$ws.Range("B381").Value = 14887; $ws.Range("C381").Value = 14749; $ws.Range("D381").Value = 15059; $ws.Range("F381").Value = 13422
$ws.Range("B382").Value = 15023; $ws.Range("C382").Value = 14862; $ws.Range("D382").Value = 15228; $ws.Range("E382").Value = 14356; $ws.Range("F382").Value = 14214; $ws.Range("G382").Value = 14524
$ws.Range("B383").Value = 14266; $ws.Range("C383").Value = 14102; $ws.Range("D383").Value = 14394; $ws.Range("E383").Value = 14902; $ws.Range("F383").Value = 14749; $ws.Range("G383").Value = 15069; $ws.Range("J383").Value = 1.27
$ws.Range("B384").Value = 14836; $ws.Range("C384").Value = 14684; $ws.Range("D384").Value = 14962; $ws.Range("E384").Value = 14753; $ws.Range("F384").Value = 14599; $ws.Range("G384").Value = 14911
$ws.Range("B385").Value = 14997; $ws.Range("C385").Value = 14860; $ws.Range("D385").Value = 15138; $ws.Range("E385").Value = 14780; $ws.Range("F385").Value = 14627; $ws.Range("G385").Value = 14930
$ws.Range("B386").Value = 15340; $ws.Range("C386").Value = 15186; $ws.Range("D386").Value = 15492; $ws.Range("E386").Value = 14860; $ws.Range("F386").Value = 14708; $ws.Range("G386").Value = 14996
$ws.Range("B387").Value = 18836; $ws.Range("C387").Value = 18670; $ws.Range("D387").Value = 18992; $ws.Range("E387").Value = 16002; $ws.Range("F387").Value = 15850; $ws.Range("G387").Value = 16146
$ws.Range("B388").Value = 18475; $ws.Range("C388").Value = 18247; $ws.Range("D388").Value = 18632; $ws.Range("F388").Value = 16741; $ws.Range("G388").Value = 17063; $ws.Range("L388").Value = 1.12
$ws.Range("B389").Value = 17875; $ws.Range("C389").Value = 17657; $ws.Range("D389").Value = 18036; $ws.Range("E389").Value = 17632; $ws.Range("F389").Value = 17440
$ws.Range("B390").Value = 16669; $ws.Range("C390").Value = 16484; $ws.Range("D390").Value = 16850; $ws.Range("E390").Value = 17964; $ws.Range("F390").Value = 17764; $ws.Range("G390").Value = 18127
$ws.Range("B391").Value = 16321; $ws.Range("C391").Value = 16167; $ws.Range("D391").Value = 16454; $ws.Range("E391").Value = 17335; $ws.Range("F391").Value = 17139; $ws.Range("G391").Value = 17493
$ws.Range("B392").Value = 15178; $ws.Range("C392").Value = 15017; $ws.Range("D392").Value = 15359; $ws.Range("E392").Value = 16511; $ws.Range("F392").Value = 16331; $ws.Range("G392").Value = 16675
$ws.Range("B393").Value = 14705; $ws.Range("C393").Value = 14477; $ws.Range("D393").Value = 14944; $ws.Range("E393").Value = 15718; $ws.Range("F393").Value = 15536; $ws.Range("G393").Value = 15902; $ws.Range("I393").Value = 0.88
$ws.Range("B394").Value = 18246; $ws.Range("C394").Value = 18001; $ws.Range("D394").Value = 18461; $ws.Range("E394").Value = 16113; $ws.Range("F394").Value = 15915; $ws.Range("G394").Value = 16304
$ws.Range("B395").Value = 17007; $ws.Range("C395").Value = 16778; $ws.Range("D395").Value = 17290; $ws.Range("E395").Value = 16284; $ws.Range("F395").Value = 16068; $ws.Range("G395").Value = 16513
$ws.Range("B396").Value = 15827; $ws.Range("C396").Value = 15688; $ws.Range("D396").Value = 16018; $ws.Range("E396").Value = 16446; $ws.Range("F396").Value = 16236; $ws.Range("G396").Value = 16678; $ws.Range("J396").Value = 1.01; $ws.Range("K396").Value = 0.96
$ws.Range("B397").Value = 16597; $ws.Range("C397").Value = 16339; $ws.Range("D397").Value = 16856; $ws.Range("E397").Value = 16919; $ws.Range("F397").Value = 16702; $ws.Range("G397").Value = 17156
$ws.Range("B398").Value = 15192; $ws.Range("C398").Value = 14907; $ws.Range("D398").Value = 15457; $ws.Range("E398").Value = 16156; $ws.Range("F398").Value = 15928; $ws.Range("G398").Value = 16405
$ws.Range("B399").Value = 14475; $ws.Range("C399").Value = 14259; $ws.Range("D399").Value = 14684; $ws.Range("E399").Value = 15523; $ws.Range("F399").Value = 15298; $ws.Range("G399").Value = 15754
$ws.Range("B400").Value = 14664; $ws.Range("C400").Value = 14368; $ws.Range("D400").Value = 14916; $ws.Range("F400").Value = 14968; $ws.Range("G400").Value = 15478; $ws.Range("J400").Value = 0.94; $ws.Range("M400").Value = 0.97
$ws.Range("B401").Value = 16664; $ws.Range("C401").Value = 16337; $ws.Range("D401").Value = 17064; $ws.Range("E401").Value = 15249; $ws.Range("F401").Value = 14968; $ws.Range("G401").Value = 15530
$ws.Range("B402").Value = 19100; $ws.Range("C402").Value = 18798; $ws.Range("D402").Value = 19506; $ws.Range("E402").Value = 16226; $ws.Range("F402").Value = 15941; $ws.Range("G402").Value = 16543; $ws.Range("J402").Value = 1.02
$ws.Range("B403").Value = 19474; $ws.Range("C403").Value = 18995; $ws.Range("D403").Value = 19976; $ws.Range("E403").Value = 17475; $ws.Range("F403").Value = 17125
$ws.Range("B404").Value = 19017; $ws.Range("C404").Value = 18599; $ws.Range("D404").Value = 19507; $ws.Range("E404").Value = 18564; $ws.Range("F404").Value = 18182; $ws.Range("G404").Value = 19013
$ws.Range("B405").Value = 18713; $ws.Range("C405").Value = 18224; $ws.Range("D405").Value = 19195; $ws.Range("E405").Value = 19076; $ws.Range("F405").Value = 18654; $ws.Range("G405").Value = 19546
$ws.Range("B406").Value = 18621; $ws.Range("C406").Value = 18034; $ws.Range("D406").Value = 19165; $ws.Range("E406").Value = 18956; $ws.Range("F406").Value = 18463; $ws.Range("G406").Value = 19461; $ws.Range("J406").Value = 1.18; $ws.Range("L406").Value = 1.1
$ws.Range("B407").Value = 18093; $ws.Range("C407").Value = 17430; $ws.Range("D407").Value = 18699; $ws.Range("E407").Value = 18611; $ws.Range("F407").Value = 18072; $ws.Range("G407").Value = 19142; $ws.Range("M407").Value = 1.15
$ws.Range("B408").Value = 22255; $ws.Range("C408").Value = 21374; $ws.Range("D408").Value = 23093; $ws.Range("F408").Value = 18765; $ws.Range("G408").Value = 20038; $ws.Range("L408").Value = 1.11; $ws.Range("M408").Value = 1.13
$ws.Range("B409").Value = 20805; $ws.Range("C409").Value = 19698; $ws.Range("D409").Value = 21889; $ws.Range("F409").Value = 19134; $ws.Range("G409").Value = 20711; $ws.Range("M409").Value = 1.1
$ws.Range("B410").Value = 20102; $ws.Range("C410").Value = 18832; $ws.Range("D410").Value = 21275; $ws.Range("E410").Value = 20314; $ws.Range("F410").Value = 19333; $ws.Range("G410").Value = 21239
$ws.Range("B411").Value = 19043; $ws.Range("C411").Value = 17891; $ws.Range("D411").Value = 20145; $ws.Range("F411").Value = 19449; $ws.Range("G411").Value = 21600; $ws.Range("L411").Value = 1
$ws.Range("B412").Value = 19985; $ws.Range("C412").Value = 18156; $ws.Range("D412").Value = 21433; $ws.Range("F412").Value = 18644; $ws.Range("G412").Value = 21185; $ws.Range("I412").Value = 0.99; $ws.Range("J412").Value = 1.06; $ws.Range("M412").Value = 1.04
$ws.Range("B413").Value = 19266; $ws.Range("C413").Value = 17427; $ws.Range("D413").Value = 20644; $ws.Range("F413").Value = 18077; $ws.Range("G413").Value = 20874; $ws.Range("I413").Value = 0.94; $ws.Range("J413").Value = 1.03; $ws.Range("M413").Value = 1.06
$ws.Range("B414").Value = 20180; $ws.Range("C414").Value = 17758; $ws.Range("D414").Value = 22509; $ws.Range("F414").Value = 17808; $ws.Range("G414").Value = 21183; $ws.Range("I414").Value = 0.92; $ws.Range("J414").Value = 1.02; $ws.Range("K414").Value = 1.04; $ws.Range("L414").Value = 1.01; $ws.Range("M414").Value = 1.08
$ws.Range("B415").Value = 23735; $ws.Range("C415").Value = 20330; $ws.Range("D415").Value = 26584; $ws.Range("E415").Value = 20792; $ws.Range("F415").Value = 18418; $ws.Range("G415").Value = 22792; $ws.Range("H415").Value = 1.01; $ws.Range("I415").Value = 0.96; $ws.Range("J415").Value = 1.08; $ws.Range("K415").Value = 1.04; $ws.Range("M415").Value = 1.09
$ws.Range("B416").Value = 22213; $ws.Range("C416").Value = 17395; $ws.Range("D416").Value = 26882; $ws.Range("E416").Value = 21348; $ws.Range("F416").Value = 18227; $ws.Range("G416").Value = 24155; $ws.Range("H416").Value = 1.07; $ws.Range("I416").Value = 0.98; $ws.Range("J416").Value = 1.16; $ws.Range("K416").Value = 1.04; $ws.Range("L416").Value = 0.98; $ws.Range("M416").Value = 1.09
$ws.Range("B417").Value = 20284; $ws.Range("C417").Value = 14454; $ws.Range("D417").Value = 25907; $ws.Range("E417").Value = 21603; $ws.Range("F417").Value = 17484; $ws.Range("G417").Value = 25470; $ws.Range("H417").Value = 1.1; $ws.Range("I417").Value = 1; $ws.Range("J417").Value = 1.23; $ws.Range("K417").Value = 1.01; $ws.Range("L417").Value = 0.93; $ws.Range("M417").Value = 1.08

# --- Add new row 418 for date 22.04.2021 (newest nowcast entry) ---
$ws.Range("A417:M417").Copy($ws.Range("A418:M418"))
$ws.Range("A418").Value = "22.04.2021"
$ws.Range("B418").Value = 17763
$ws.Range("C418").Value = 11684
$ws.Range("D418").Value = 24243
$ws.Range("E418").Value = 20999
$ws.Range("F418").Value = 15966
$ws.Range("G418").Value = 25904
$ws.Range("H418").Value = 1.07
$ws.Range("I418").Value = 0.91
$ws.Range("J418").Value = 1.22
$ws.Range("K418").Value = "."
$ws.Range("L418").Value = "."
$ws.Range("M418").Value = "."

# --- Update the visible selection to include the new row ---
$ws.Range("H6:M418").Select()
